# Trade #3 closed at 2026-02-16 22:56:37 - base_strategy UP +0.000%
# Append the new trade row (row 4) to both the "All Trades" sheet and the
# per-strategy "base_strategy" sheet - they mirror each other's trade log.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A4").Value = 3
    $ws.Range("B4").Value = "'2026-02-16"
    $ws.Range("C4").Value = "'22:56:37"
    $ws.Range("D4").Value = "base_strategy"
    $ws.Range("E4").Value = "UP"
    $ws.Range("F4").Value = 0.5
    $ws.Range("G4").Value = "'"
    $ws.Range("H4").Value = "OPEN"
    $ws.Range("I4").Value = 0
    $ws.Range("J4").Value = 0
    $ws.Range("K4").Value = 100
    $ws.Range("L4").Value = 0
    $ws.Range("M4").Value = 0
    $ws.Range("N4").Value = 0.6
    $ws.Range("O4").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P4").Value = "'"
    $ws.Range("Q4").Value = 0
}
